$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (was: Dina / Primera / 44161) -> becomes Castle Brite / Primera / 44160
$ws.Range("D4").Value = 44160
$ws.Range("K4").Value = "Castle Brite"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 240
$ws.Range("N4").Value = 20500
$ws.Range("O4").Value = 21000
$ws.Range("P4").Value = 20750
$ws.Range("Q4").Value = "$/caja 15 kilos"
$ws.Range("S4").Value = 1383
$ws.Range("T4").Value = 15

# Row 5 (was: Dina / Segunda / 44161) -> becomes Castle Brite / Primera / 44175
$ws.Range("D5").Value = 44175
$ws.Range("K5").Value = "Castle Brite"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 300
$ws.Range("N5").Value = 21000
$ws.Range("O5").Value = 22000
$ws.Range("P5").Value = 21500
$ws.Range("Q5").Value = "$/caja 18 kilos"
$ws.Range("S5").Value = 1194
$ws.Range("T5").Value = 18

# Row 6 (was: Castle Brite / Primera / 44160) -> becomes Dina / Primera / 44161
$ws.Range("D6").Value = 44161
$ws.Range("K6").Value = "Dina"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 300
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 20500
$ws.Range("P6").Value = 20250
$ws.Range("Q6").Value = "$/caja 15 kilos"
$ws.Range("S6").Value = 1350
$ws.Range("T6").Value = 15

# Row 7 (was: Castle Brite / Primera / 44175) -> becomes Dina / Segunda / 44161
$ws.Range("D7").Value = 44161
$ws.Range("K7").Value = "Dina"
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 18000
$ws.Range("O7").Value = 18500
$ws.Range("P7").Value = 18250
$ws.Range("Q7").Value = "$/caja 15 kilos"
$ws.Range("S7").Value = 1217
$ws.Range("T7").Value = 15
